{"js": "// Rewrite the Summary, Experience, Education, Skills and Projects body\n// paragraphs of the resume with new copy. Each resume section is a\n// Heading1 paragraph (\"Summary\", \"Experience\", ...) immediately followed\n// by a single body paragraph holding the section's text; that body\n// paragraph is located via getNext() off the heading and is fully\n// replaced with Word.InsertLocation.replace so the resulting markup is a\n// single run (line-break characters \\u000b become <w:br/> elements).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Build a lookup: heading text (trimmed) -> paragraph object immediately\n// following it (the section's content paragraph).\nconst sectionBodyParagraph = {};\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length - 1; i++) {\n  const text = items[i].text.trim();\n  if (text === \"Summary\" || text === \"Experience\" || text === \"Education\" ||\n      text === \"Skills\" || text === \"Projects\") {\n    sectionBodyParagraph[text] = items[i + 1];\n  }\n}\n\nconst LB = \"\\u000b\"; // manual line break -> <w:br/>\n\nconst newSummary =\n  \"Highly motivated and results-oriented individual seeking a challenging \" +\n  \"role in the development of innovative mobile applications. Proven \" +\n  \"ability to collaborate effectively, solve complex problems, and \" +\n  \"contribute to a dynamic team environment. Eager to leverage skills in \" +\n  \"software development, project management, and communication to \" +\n  \"contribute to impactful projects.\";\n\nconst newExperience = [\n  \"AI Intern, Blue Silicon Infotech\",\n  \"Developed and implemented AI resume enhancer, resulting in a 20% \" +\n    \"increase in resume completion rates. Optimized resume templates for \" +\n    \"improved readability and clarity. Achieved quantifiable results, \" +\n    \"including a 15% improvement in resume accuracy.\",\n].join(LB);\n\nconst newEducation = \"Bachelor of Engineering from AVIT. Graduated: 2026-05. GPA: 7.1.\";\n\nconst newSkills =\n  \"Here's a revised skills section tailored for a global editing resume, \" +\n  \"focusing on conciseness, organization, and professionalism:, Skills**, \" +\n  \"   **Programming Languages:** Python, Java, Node.js, TypeScript, C#, \" +\n  \"Flutter, Kotlin, Dart, React, Python, SQL,    **Database:** MySQL, \" +\n  \"PostgreSQL, MongoDB, SQL Server,    **Web Development:** HTML, CSS, \" +\n  \"JavaScript, React, Angular, Vue.js,    **Cloud Technologies:** AWS, \" +\n  \"Azure, Google Cloud Platform,    **Operating Systems:** Linux, \" +\n  \"Windows, macOS,    **Version Control:** Git, GitHub, GitLab,    \" +\n  \"**Data Analysis:** Pandas, NumPy, Matplotlib, Seaborn,    **Testing:** \" +\n  \"Unit Testing, Integration Testing, End-to-End Testing,    **API \" +\n  \"Development:** RESTful APIs, GraphQL,    **Design Principles:** SOLID, \" +\n  \"DRY, KISS,    **Other:** Agile Development, Mobile Development, Data \" +\n  \"Science\";\n\nconst newProjects = [\n  \"**Project:** Enhanced QR Scanner and Generator\",\n  \"**Summary:** This project aimed to significantly improve the \" +\n    \"efficiency and accuracy of QR scanning and generator functionality. \" +\n    \"By implementing a novel algorithm and incorporating real-time data \" +\n    \"integration, we achieved a demonstrable increase in accuracy and \" +\n    \"reduced processing time. This improved functionality was directly \" +\n    \"translated into increased sales and reduced operational costs.\",\n  \"**Technologies:**\",\n  \"*  QR scanner and generator\",\n  \"*  Prediction pro\",\n  \"*  Simple purchase order manager\",\n  \"*  PDF maker\",\n  \"**Contributions:**\",\n  \"*  Improved accuracy in QR scanning and generator processing.\",\n  \"*  Enhanced real-time data integration for improved processing.\",\n  \"*  Reduced processing time by 20%.\",\n  \"**Measurable Results:**\",\n  \"*  Increased accuracy in QR scanning and generator processing.\",\n  \"*  Reduced processing time by 20%.\",\n  \"*  Improved sales and reduced operational costs.\",\n  \"**Improvements:**\",\n  \"*  Improved accuracy in QR scanning and generator processing.\",\n  \"*  Enhanced real-time data integration for improved processing.\",\n  \"*  Reduced processing time by 20%.\",\n].join(LB);\n\nsectionBodyParagraph[\"Summary\"].insertText(newSummary, Word.InsertLocation.replace);\nsectionBodyParagraph[\"Experience\"].insertText(newExperience, Word.InsertLocation.replace);\nsectionBodyParagraph[\"Education\"].insertText(newEducation, Word.InsertLocation.replace);\nsectionBodyParagraph[\"Skills\"].insertText(newSkills, Word.InsertLocation.replace);\nsectionBodyParagraph[\"Projects\"].insertText(newProjects, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Rewrite the Summary, Experience, Education, Skills and Projects body\n# paragraphs of the resume with new copy. Each resume section is a\n# \"Heading 1\" paragraph (\"Summary\", \"Experience\", ...) immediately\n# followed by a single body paragraph holding the section's text; that\n# body paragraph is located by scanning $doc.Paragraphs for the heading\n# text and taking the next paragraph. Assigning .Range.Text directly\n# replaces the whole paragraph's run content in one shot, and embedding\n# [char]11 (vertical tab, Word's manual line-break character) produces\n# <w:br/> elements instead of new paragraphs.\n\n$d = $word.ActiveDocument\n$LB = [char]11\n\n$count = $d.Paragraphs.Count\n$sectionBodyIndex = @{}\nfor ($i = 1; $i -lt $count; $i++) {\n    $headingText = $d.Paragraphs($i).Range.Text.Trim()\n    if ($headingText -eq \"Summary\" -or $headingText -eq \"Experience\" -or `\n        $headingText -eq \"Education\" -or $headingText -eq \"Skills\" -or `\n        $headingText -eq \"Projects\") {\n        $sectionBodyIndex[$headingText] = $i + 1\n    }\n}\n\n$newSummary = \"Highly motivated and results-oriented individual seeking a challenging \" + `\n    \"role in the development of innovative mobile applications. Proven \" + `\n    \"ability to collaborate effectively, solve complex problems, and \" + `\n    \"contribute to a dynamic team environment. Eager to leverage skills in \" + `\n    \"software development, project management, and communication to \" + `\n    \"contribute to impactful projects.\"\n\n$newExperience = \"AI Intern, Blue Silicon Infotech\" + $LB + `\n    \"Developed and implemented AI resume enhancer, resulting in a 20% \" + `\n    \"increase in resume completion rates. Optimized resume templates for \" + `\n    \"improved readability and clarity. Achieved quantifiable results, \" + `\n    \"including a 15% improvement in resume accuracy.\"\n\n$newEducation = \"Bachelor of Engineering from AVIT. Graduated: 2026-05. GPA: 7.1.\"\n\n$newSkills = \"Here's a revised skills section tailored for a global editing resume, \" + `\n    \"focusing on conciseness, organization, and professionalism:, Skills**, \" + `\n    \"   **Programming Languages:** Python, Java, Node.js, TypeScript, C#, \" + `\n    \"Flutter, Kotlin, Dart, React, Python, SQL,    **Database:** MySQL, \" + `\n    \"PostgreSQL, MongoDB, SQL Server,    **Web Development:** HTML, CSS, \" + `\n    \"JavaScript, React, Angular, Vue.js,    **Cloud Technologies:** AWS, \" + `\n    \"Azure, Google Cloud Platform,    **Operating Systems:** Linux, \" + `\n    \"Windows, macOS,    **Version Control:** Git, GitHub, GitLab,    \" + `\n    \"**Data Analysis:** Pandas, NumPy, Matplotlib, Seaborn,    **Testing:** \" + `\n    \"Unit Testing, Integration Testing, End-to-End Testing,    **API \" + `\n    \"Development:** RESTful APIs, GraphQL,    **Design Principles:** SOLID, \" + `\n    \"DRY, KISS,    **Other:** Agile Development, Mobile Development, Data \" + `\n    \"Science\"\n\n$newProjects = @(\n    \"**Project:** Enhanced QR Scanner and Generator\",\n    (\"**Summary:** This project aimed to significantly improve the \" + `\n        \"efficiency and accuracy of QR scanning and generator functionality. \" + `\n        \"By implementing a novel algorithm and incorporating real-time data \" + `\n        \"integration, we achieved a demonstrable increase in accuracy and \" + `\n        \"reduced processing time. This improved functionality was directly \" + `\n        \"translated into increased sales and reduced operational costs.\"),\n    \"**Technologies:**\",\n    \"*  QR scanner and generator\",\n    \"*  Prediction pro\",\n    \"*  Simple purchase order manager\",\n    \"*  PDF maker\",\n    \"**Contributions:**\",\n    \"*  Improved accuracy in QR scanning and generator processing.\",\n    \"*  Enhanced real-time data integration for improved processing.\",\n    \"*  Reduced processing time by 20%.\",\n    \"**Measurable Results:**\",\n    \"*  Increased accuracy in QR scanning and generator processing.\",\n    \"*  Reduced processing time by 20%.\",\n    \"*  Improved sales and reduced operational costs.\",\n    \"**Improvements:**\",\n    \"*  Improved accuracy in QR scanning and generator processing.\",\n    \"*  Enhanced real-time data integration for improved processing.\",\n    \"*  Reduced processing time by 20%.\"\n) -join $LB\n\n$d.Paragraphs($sectionBodyIndex[\"Summary\"]).Range.Text = $newSummary\n$d.Paragraphs($sectionBodyIndex[\"Experience\"]).Range.Text = $newExperience\n$d.Paragraphs($sectionBodyIndex[\"Education\"]).Range.Text = $newEducation\n$d.Paragraphs($sectionBodyIndex[\"Skills\"]).Range.Text = $newSkills\n$d.Paragraphs($sectionBodyIndex[\"Projects\"]).Range.Text = $newProjects\n"}
